$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 702.5454999999999
$ws.Range("I2").Value = 672.8
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 672.8
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -559.8
$ws.Range("N2").Value = -1226
$ws.Range("H29").Value = 10241.917
$ws.Range("I29").Value = 1116.5
$ws.Range("J29").Value = 19367.334
$ws.Range("K29").Value = 3349.5
$ws.Range("L29").Value = 58102.00199999999
$ws.Range("M29").Value = -3068.5
$ws.Range("N29").Value = -58664.00199999999
$ws.Range("H38").Value = 2432.2727
$ws.Range("I38").Value = 149
$ws.Range("J38").Value = 3737
$ws.Range("K38").Value = 447
$ws.Range("L38").Value = 11211
$ws.Range("M38").Value = -75
$ws.Range("N38").Value = -11955
$ws.Range("H41").Value = 765.06665
$ws.Range("I41").Value = 536.3333
$ws.Range("J41").Value = 1108.1666
$ws.Range("K41").Value = 536.3333
$ws.Range("L41").Value = 1108.1666
$ws.Range("M41").Value = -96.33330000000001
$ws.Range("N41").Value = -1988.1666
$ws.Range("H58").Value = 2026.0769
$ws.Range("I58").Value = 1337.8
$ws.Range("J58").Value = 2456.25
$ws.Range("K58").Value = 4013.4
$ws.Range("L58").Value = 7368.75
$ws.Range("M58").Value = -3863.4
$ws.Range("N58").Value = -7668.75
$ws.Range("H64").Value = 5085.905
$ws.Range("I64").Value = 5495.7144
$ws.Range("J64").Value = 4266.2856
$ws.Range("K64").Value = 5495.7144
$ws.Range("L64").Value = 4266.2856
$ws.Range("M64").Value = -5247.7144
$ws.Range("N64").Value = -4762.2856
$ws.Range("H67").Value = 5085.905
$ws.Range("I67").Value = 5495.7144
$ws.Range("J67").Value = 4266.2856
$ws.Range("K67").Value = 5495.7144
$ws.Range("L67").Value = 4266.2856
$ws.Range("M67").Value = -4637.7144
$ws.Range("N67").Value = -5982.2856
$ws.Range("H87").Value = 32000
$ws.Range("J87").Value = 32000
$ws.Range("L87").Value = 32000
$ws.Range("N87").Value = -34496
$ws.Range("H90").Value = 32000
$ws.Range("J90").Value = 32000
$ws.Range("L90").Value = 96000
$ws.Range("N90").Value = -108480
$ws.Range("H98").Value = 1751.25
$ws.Range("I98").Value = 1751.25
$ws.Range("K98").Value = 1751.25
$ws.Range("M98").Value = -253.25
$ws.Range("H112").Value = 901.35187
$ws.Range("J112").Value = 935.2708
$ws.Range("L112").Value = 2805.8124
$ws.Range("N112").Value = -5021.8124
$ws.Range("H122").Value = 1751.25
$ws.Range("I122").Value = 1751.25
$ws.Range("K122").Value = 5253.75
$ws.Range("M122").Value = -2803.75
$ws.Range("H132").Value = 1195.1082
$ws.Range("I132").Value = 1249.0883
$ws.Range("J132").Value = 583.3333
$ws.Range("K132").Value = 3747.2649
$ws.Range("L132").Value = 1749.9999
$ws.Range("M132").Value = -1217.2649
$ws.Range("N132").Value = -6809.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3851.959
$ws.Range("I32").Value = 3240.1667
$ws.Range("J32").Value = 5590.737
$ws.Range("K32").Value = 3240.1667
$ws.Range("L32").Value = 5590.737
$ws.Range("M32").Value = -2953.1667
$ws.Range("N32").Value = -6164.737
$ws.Range("H106").Value = 49444
$ws.Range("J106").Value = 49444
$ws.Range("L106").Value = 49444
$ws.Range("N106").Value = -51968
$ws.Range("H110").Value = 1289
$ws.Range("I110").Value = 1202.2
$ws.Range("J110").Value = 1506
$ws.Range("K110").Value = 1202.2
$ws.Range("L110").Value = 1506
$ws.Range("M110").Value = 842.8
$ws.Range("N110").Value = -5596

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 125002024
$ws.Range("I20").Value = 333334500
$ws.Range("K20").Value = 333334500
$ws.Range("M20").Value = -333334253
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H109").Value = 54000
$ws.Range("J109").Value = 54000
$ws.Range("L109").Value = 54000
$ws.Range("N109").Value = -56774
$ws.Range("H138").Value = 59733.332
$ws.Range("J138").Value = 59733.332
$ws.Range("L138").Value = 59733.332
$ws.Range("N138").Value = -70013.33199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2587.6301
$ws.Range("I31").Value = 2225.7058
$ws.Range("J31").Value = 2697.5
$ws.Range("K31").Value = 2225.7058
$ws.Range("L31").Value = 2697.5
$ws.Range("M31").Value = -1930.7058
$ws.Range("N31").Value = -3287.5
$ws.Range("H34").Value = 2587.6301
$ws.Range("I34").Value = 2225.7058
$ws.Range("J34").Value = 2697.5
$ws.Range("K34").Value = 2225.7058
$ws.Range("L34").Value = 2697.5
$ws.Range("M34").Value = -2023.7058
$ws.Range("N34").Value = -3101.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3573.3962
$ws.Range("I68").Value = 4988.48
$ws.Range("J68").Value = 2309.9285
$ws.Range("K68").Value = 14965.44
$ws.Range("L68").Value = 6929.7855
$ws.Range("M68").Value = -14154.44
$ws.Range("N68").Value = -8551.7855
$ws.Range("H71").Value = 3573.3962
$ws.Range("I71").Value = 4988.48
$ws.Range("J71").Value = 2309.9285
$ws.Range("K71").Value = 44896.31999999999
$ws.Range("L71").Value = 20789.3565
$ws.Range("M71").Value = -40840.31999999999
$ws.Range("N71").Value = -28901.3565
$ws.Range("H80").Value = 6200
$ws.Range("H83").Value = 6200
$ws.Range("H107").Value = 791.2
$ws.Range("I107").Value = 279.7143
$ws.Range("J107").Value = 1035.3182
$ws.Range("K107").Value = 839.1428999999999
$ws.Range("L107").Value = 3105.9546
$ws.Range("M107").Value = 1080.8571
$ws.Range("N107").Value = -6945.9546
$ws.Range("H112").Value = 2108.5
$ws.Range("J112").Value = 3800
$ws.Range("L112").Value = 11400
$ws.Range("N112").Value = -13616
$ws.Range("H131").Value = 18645244
$ws.Range("I131").Value = 6667086.5
$ws.Range("J131").Value = 22728708
$ws.Range("K131").Value = 20001259.5
$ws.Range("L131").Value = 68186124
$ws.Range("M131").Value = -19996219.5
$ws.Range("N131").Value = -68196204
$ws.Range("H137").Value = 30316970
$ws.Range("I137").Value = 1167
$ws.Range("J137").Value = 50022240
$ws.Range("K137").Value = 3501
$ws.Range("L137").Value = 150066720
$ws.Range("M137").Value = 1599
$ws.Range("N137").Value = -150076920

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2892.6155
$ws.Range("I7").Value = 2434
$ws.Range("J7").Value = 3285.7144
$ws.Range("K7").Value = 2434
$ws.Range("L7").Value = 3285.7144
$ws.Range("M7").Value = -2322
$ws.Range("N7").Value = -3509.7144
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H68").Value = 55558316
$ws.Range("I68").Value = 2890
$ws.Range("J68").Value = 166669170
$ws.Range("K68").Value = 2890
$ws.Range("L68").Value = 166669170
$ws.Range("M68").Value = -2141
$ws.Range("N68").Value = -166670668
$ws.Range("H71").Value = 55558316
$ws.Range("I71").Value = 2890
$ws.Range("J71").Value = 166669170
$ws.Range("K71").Value = 14450
$ws.Range("L71").Value = 833345850
$ws.Range("M71").Value = -10706
$ws.Range("N71").Value = -833353338
$ws.Range("H93").Value = 100041800
$ws.Range("I93").Value = 68001
$ws.Range("K93").Value = 68001
$ws.Range("M93").Value = -66753
$ws.Range("H126").Value = 2892.6155
$ws.Range("I126").Value = 2434
$ws.Range("J126").Value = 3285.7144
$ws.Range("K126").Value = 7302
$ws.Range("L126").Value = 9857.143199999999
$ws.Range("M126").Value = -4832
$ws.Range("N126").Value = -14797.1432
$ws.Range("H136").Value = 4517.6665
$ws.Range("I136").Value = 2399.1035
$ws.Range("K136").Value = 7197.310500000001
$ws.Range("M136").Value = -4647.310500000001
